# Updates the cryptos list (price / 1h volume columns, plus the
# Monero <-> dogwifhat row swap at rows 41/42) to match the refreshed
# data, mirroring the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '62.893.08'
$ws.Cells.Item(2, 5).Value = '  +0.29%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.462.42'
$ws.Cells.Item(3, 5).Value = '  +0.65%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '574.00'
$ws.Cells.Item(5, 5).Value = '  -0.71%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '146.43'
$ws.Cells.Item(6, 5).Value = '  +0.47%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.04%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.17%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '2.462.86'
$ws.Cells.Item(9, 5).Value = '  +0.73%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.111'
$ws.Cells.Item(10, 5).Value = '  +0.67%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.86%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.69%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +1.06%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '28.96'
$ws.Cells.Item(14, 5).Value = '  +1.97%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -0.65%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '2.909.37'
$ws.Cells.Item(16, 5).Value = '  +0.66%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '62.780.96'
$ws.Cells.Item(17, 5).Value = '  +0.32%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '2.465.66'
$ws.Cells.Item(18, 5).Value = '  +0.77%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '7.99'
$ws.Cells.Item(19, 5).Value = '  +2.21%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.73%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '326.93'
$ws.Cells.Item(21, 5).Value = '  -0.50%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '2.22'
$ws.Cells.Item(22, 5).Value = '  +9.88%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -0.04%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '10.11'
$ws.Cells.Item(25, 5).Value = '  +19.41%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '65.72'
$ws.Cells.Item(26, 5).Value = '  +0.20%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '649.91'
$ws.Cells.Item(27, 5).Value = '  +0.98%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.0₃0983'
$ws.Cells.Item(28, 5).Value = '  -0.43%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.584.51'
$ws.Cells.Item(29, 5).Value = '  +0.42%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -14.34%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.54%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '7.99'
$ws.Cells.Item(32, 5).Value = '  -2.56%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -1.29%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -4.03%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -0.01%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +3.67%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '4.75'
$ws.Cells.Item(37, 5).Value = '  +0.04%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '5.42'
$ws.Cells.Item(38, 5).Value = '  -1.23%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -1.28%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '18.72'
$ws.Cells.Item(40, 5).Value = '  +0.41%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Monero'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '150.36'
$ws.Cells.Item(41, 5).Value = '  -1.97%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'dogwifhat'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '2.77'
$ws.Cells.Item(42, 5).Value = '  +1.52%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.74'
$ws.Cells.Item(43, 5).Value = '  -1.28%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.0₆0309'
$ws.Cells.Item(44, 5).Value = '  -43.71%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.01%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '153.49'
$ws.Cells.Item(46, 5).Value = '  +5.74%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '15.24'

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '3.59'
$ws.Cells.Item(48, 5).Value = '  -0.68%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '20.54'
$ws.Cells.Item(49, 5).Value = '  -0.91%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.607'
$ws.Cells.Item(50, 5).Value = '  +0.38%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.81%  '
